{"js": "// Update the worksheet date and all 25 two-digit-by-two-digit\n// multiplication prompts to the values generated for the new date.\n// Each old text is unique in the document, so a simple exact-match\n// search + full-text replace (preserving the run's formatting) is\n// applied for every pair.\n\nconst replacements = [\n  [\"2024-12-19 Thursday\", \"2024-12-20 Friday\"],\n  [\"90\u00d748=\", \"40\u00d756=\"],\n  [\"23\u00d764=\", \"55\u00d714=\"],\n  [\"46\u00d791=\", \"46\u00d726=\"],\n  [\"24\u00d730=\", \"78\u00d742=\"],\n  [\"95\u00d750=\", \"11\u00d745=\"],\n  [\"59\u00d743=\", \"39\u00d791=\"],\n  [\"70\u00d743=\", \"78\u00d715=\"],\n  [\"26\u00d799=\", \"36\u00d772=\"],\n  [\"63\u00d729=\", \"45\u00d738=\"],\n  [\"41\u00d791=\", \"46\u00d786=\"],\n  [\"81\u00d752=\", \"14\u00d770=\"],\n  [\"28\u00d775=\", \"17\u00d770=\"],\n  [\"97\u00d772=\", \"22\u00d784=\"],\n  [\"51\u00d781=\", \"46\u00d793=\"],\n  [\"44\u00d793=\", \"83\u00d759=\"],\n  [\"92\u00d732=\", \"56\u00d713=\"],\n  [\"81\u00d751=\", \"65\u00d718=\"],\n  [\"71\u00d783=\", \"91\u00d731=\"],\n  [\"46\u00d790=\", \"21\u00d787=\"],\n  [\"82\u00d771=\", \"48\u00d730=\"],\n  [\"70\u00d768=\", \"90\u00d713=\"],\n  [\"28\u00d725=\", \"12\u00d740=\"],\n  [\"68\u00d721=\", \"90\u00d797=\"],\n  [\"39\u00d716=\", \"73\u00d731=\"],\n  [\"55\u00d793=\", \"25\u00d787=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 two-digit-by-two-digit\n# multiplication prompts to the values generated for the new date.\n# Each old text is unique in the document, so a simple Find/Replace\n# (wdReplaceAll) is used for every pair; formatting of the existing\n# run is preserved by Find/Replace.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-12-19 Thursday\", \"2024-12-20 Friday\"),\n  @(\"90\u00d748=\", \"40\u00d756=\"),\n  @(\"23\u00d764=\", \"55\u00d714=\"),\n  @(\"46\u00d791=\", \"46\u00d726=\"),\n  @(\"24\u00d730=\", \"78\u00d742=\"),\n  @(\"95\u00d750=\", \"11\u00d745=\"),\n  @(\"59\u00d743=\", \"39\u00d791=\"),\n  @(\"70\u00d743=\", \"78\u00d715=\"),\n  @(\"26\u00d799=\", \"36\u00d772=\"),\n  @(\"63\u00d729=\", \"45\u00d738=\"),\n  @(\"41\u00d791=\", \"46\u00d786=\"),\n  @(\"81\u00d752=\", \"14\u00d770=\"),\n  @(\"28\u00d775=\", \"17\u00d770=\"),\n  @(\"97\u00d772=\", \"22\u00d784=\"),\n  @(\"51\u00d781=\", \"46\u00d793=\"),\n  @(\"44\u00d793=\", \"83\u00d759=\"),\n  @(\"92\u00d732=\", \"56\u00d713=\"),\n  @(\"81\u00d751=\", \"65\u00d718=\"),\n  @(\"71\u00d783=\", \"91\u00d731=\"),\n  @(\"46\u00d790=\", \"21\u00d787=\"),\n  @(\"82\u00d771=\", \"48\u00d730=\"),\n  @(\"70\u00d768=\", \"90\u00d713=\"),\n  @(\"28\u00d725=\", \"12\u00d740=\"),\n  @(\"68\u00d721=\", \"90\u00d797=\"),\n  @(\"39\u00d716=\", \"73\u00d731=\"),\n  @(\"55\u00d793=\", \"25\u00d787=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $r = $d.Content\n  $r.Find.ClearFormatting()\n  $r.Find.Replacement.ClearFormatting()\n  $r.Find.Text = $oldText\n  $r.Find.Replacement.Text = $newText\n  $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
